$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "29.943.98"
Set-TextValue "D3" "1.633.14"
Set-TextValue "D5" "214.90"
Set-TextValue "E5" "  +1.18%  "
Set-TextValue "E6" "  +1.27%  "
Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  -0.09%  "
Set-TextValue "D8" "29.61"
Set-TextValue "E8" "  +10.03%  "
Set-TextValue "E9" "  +3.68%  "
Set-TextValue "D10" "0.0614"
Set-TextValue "E11" "  +0.69%  "
Set-TextValue "D12" "1.865.71"
Set-TextValue "E12" "  +2.05%  "
Set-TextValue "D13" "1.631.61"
Set-TextValue "E13" "  +1.95%  "
Set-TextValue "E14" "  +6.52%  "
Set-TextValue "D15" "9.45"
Set-TextValue "E15" "  +24.01%  "
Set-TextValue "D16" "3.90"
Set-TextValue "E16" "  +4.49%  "
Set-TextValue "D17" "29.953.77"
Set-TextValue "E17" "  +1.07%  "
Set-TextValue "D18" "64.87"
Set-TextValue "E18" "  +1.42%  "
Set-TextValue "D19" "249.68"
Set-TextValue "E19" "  +3.31%  "
Set-TextValue "D20" "0.0₃0707"
Set-TextValue "E20" "  +1.94%  "
Set-TextValue "E21" "  -0.09%  "
Set-TextValue "D22" "4.17"
Set-TextValue "E22" "  +5.02%  "
Set-TextValue "D23" "9.64"
Set-TextValue "E23" "  +4.28%  "
Set-TextValue "D25" "159.75"
Set-TextValue "E25" "  +3.35%  "
Set-TextValue "E26" "  +2.44%  "
Set-TextValue "E27" "  +2.53%  "
Set-TextValue "D28" "6.62"
Set-TextValue "E28" "  +3.57%  "
Set-TextValue "E30" "  +2.67%  "
Set-TextValue "E31" "  +6.56%  "
Set-TextValue "D32" "3.38"
Set-TextValue "E32" "  +4.82%  "
Set-TextValue "E33" "  +2.15%  "
Set-TextValue "D34" "1.430.88"
Set-TextValue "E34" "  +0.50%  "
Set-TextValue "E35" "  +7.29%  "
Set-TextValue "E36" "  +1.58%  "
Set-TextValue "E37" "  -0.34%  "
Set-TextValue "E38" "  +3.31%  "
Set-TextValue "E39" "  -0.27%  "
Set-TextValue "D40" "0.557"
Set-TextValue "E40" "  +2.51%  "
Set-TextValue "D41" "73.38"
Set-TextValue "E41" "  +11.26%  "
Set-TextValue "D42" "1.98"
Set-TextValue "E42" "  +1.04%  "
Set-TextValue "D43" "0.830"
Set-TextValue "E43" "  +3.10%  "
Set-TextValue "D44" "54.76"
Set-TextValue "E44" "  +0.93%  "
Set-TextValue "D45" "0.0495"
Set-TextValue "E45" "  +0.16%  "
Set-TextValue "E46" "  +5.46%  "
Set-TextValue "E47" "  -0.09%  "
Set-TextValue "E48" "  +3.22%  "
Set-TextValue "D49" "1.771.67"
Set-TextValue "E49" "  +1.20%  "
Set-TextValue "D50" "89.95"
Set-TextValue "E50" "  +4.52%  "
Set-TextValue "E51" "  +3.54%  "
